# Auto-generated edit script: update cryptos list values
# Applies text-preserving writes (leading apostrophe + style reset)
# so numeric-looking strings (e.g. '58.436.79', '0.0000166', '  +1.28%  ')
# remain literal text instead of being reinterpreted as numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.Value = '''58.436.79'
$cell.Style = 'Normal'
$cell = $ws.Range('E2')
$cell.Value = '''  -0.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D3')
$cell.Value = '''3.137.32'
$cell.Style = 'Normal'
$cell = $ws.Range('E3')
$cell.Value = '''  +1.29%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E4')
$cell.Value = '''  +0.01%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D5')
$cell.Value = '''533.90'
$cell.Style = 'Normal'
$cell = $ws.Range('E5')
$cell.Value = '''  +1.28%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D6')
$cell.Value = '''143.06'
$cell.Style = 'Normal'
$cell = $ws.Range('E6')
$cell.Value = '''  +0.31%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D7')
$cell.Value = '''1.00'
$cell.Style = 'Normal'
$cell = $ws.Range('E7')
$cell.Value = '''  +0.00%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D8')
$cell.Value = '''3.139.33'
$cell.Style = 'Normal'
$cell = $ws.Range('E8')
$cell.Value = '''  +1.39%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D9')
$cell.Value = '''0.451'
$cell.Style = 'Normal'
$cell = $ws.Range('E9')
$cell.Value = '''  +2.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D10')
$cell.Value = '''7.17'
$cell.Style = 'Normal'
$cell = $ws.Range('E10')
$cell.Value = '''  -2.18%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E11')
$cell.Value = '''  +0.71%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D12')
$cell.Value = '''0.394'
$cell.Style = 'Normal'
$cell = $ws.Range('E12')
$cell.Value = '''  +3.01%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D13')
$cell.Value = '''3.680.61'
$cell.Style = 'Normal'
$cell = $ws.Range('E13')
$cell.Value = '''  +1.36%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E14')
$cell.Value = '''  +3.34%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D15')
$cell.Value = '''25.65'
$cell.Style = 'Normal'
$cell = $ws.Range('E15')
$cell.Value = '''  -4.47%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D16')
$cell.Value = '''0.0000166'
$cell.Style = 'Normal'
$cell = $ws.Range('E16')
$cell.Value = '''  +0.23%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D17')
$cell.Value = '''58.518.59'
$cell.Style = 'Normal'
$cell = $ws.Range('E17')
$cell.Value = '''  -0.07%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D18')
$cell.Value = '''3.140.38'
$cell.Style = 'Normal'
$cell = $ws.Range('E18')
$cell.Value = '''  +1.37%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E19')
$cell.Value = '''  +0.26%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D20')
$cell.Value = '''12.89'
$cell.Style = 'Normal'
$cell = $ws.Range('E20')
$cell.Value = '''  -0.22%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D21')
$cell.Value = '''7.99'
$cell.Style = 'Normal'
$cell = $ws.Range('E21')
$cell.Value = '''  -0.81%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D22')
$cell.Value = '''342.66'
$cell.Style = 'Normal'
$cell = $ws.Range('E22')
$cell.Value = '''  +0.41%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E23')
$cell.Value = '''  +0.24%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D24')
$cell.Value = '''0.513'
$cell.Style = 'Normal'
$cell = $ws.Range('E24')
$cell.Value = '''  +1.87%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D25')
$cell.Value = '''67.79'
$cell.Style = 'Normal'
$cell = $ws.Range('E25')
$cell.Value = '''  +2.60%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E26')
$cell.Value = '''  -0.42%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E27')
$cell.Value = '''  -0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D28')
$cell.Value = '''0.0₃0931'
$cell.Style = 'Normal'
$cell = $ws.Range('E28')
$cell.Value = '''  +2.00%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D29')
$cell.Value = '''7.51'
$cell.Style = 'Normal'
$cell = $ws.Range('E29')
$cell.Value = '''  +3.57%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E30')
$cell.Value = '''  -2.45%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E31')
$cell.Value = '''  +0.04%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E32')
$cell.Value = '''  +1.40%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E33')
$cell.Value = '''  +1.15%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E34')
$cell.Value = '''  -0.24%  '
$cell.Style = 'Normal'
$cell = $ws.Range('B35')
$cell.Value = '''Monero'
$cell.Style = 'Normal'
$cell = $ws.Range('C35')
$cell.Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.Style = 'Normal'
$cell = $ws.Range('D35')
$cell.Value = '''158.26'
$cell.Style = 'Normal'
$cell = $ws.Range('E35')
$cell.Value = '''  +2.63%  '
$cell.Style = 'Normal'
$cell = $ws.Range('B36')
$cell.Value = '''NEARProtocol'
$cell.Style = 'Normal'
$cell = $ws.Range('C36')
$cell.Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.Style = 'Normal'
$cell = $ws.Range('D36')
$cell.Value = '''4.80'
$cell.Style = 'Normal'
$cell = $ws.Range('E36')
$cell.Value = '''  +3.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D37')
$cell.Value = '''6.26'
$cell.Style = 'Normal'
$cell = $ws.Range('E37')
$cell.Value = '''  +3.41%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D38')
$cell.Value = '''26.18'
$cell.Style = 'Normal'
$cell = $ws.Range('E38')
$cell.Value = '''  -2.76%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E39')
$cell.Value = '''  -3.45%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E40')
$cell.Value = '''  +11.63%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D41')
$cell.Value = '''0.0672'
$cell.Style = 'Normal'
$cell = $ws.Range('E41')
$cell.Value = '''  -0.74%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D42')
$cell.Value = '''0.709'
$cell.Style = 'Normal'
$cell = $ws.Range('E42')
$cell.Value = '''  +4.66%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E43')
$cell.Value = '''  +3.58%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.Value = '''3.180.79'
$cell.Style = 'Normal'
$cell = $ws.Range('E44')
$cell.Value = '''  +1.26%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.Value = '''36.61'
$cell.Style = 'Normal'
$cell = $ws.Range('E45')
$cell.Value = '''  -0.41%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E46')
$cell.Value = '''  +0.06%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.Value = '''0.0265'
$cell.Style = 'Normal'
$cell = $ws.Range('E47')
$cell.Value = '''  +3.25%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D48')
$cell.Value = '''2.302.39'
$cell.Style = 'Normal'
$cell = $ws.Range('E48')
$cell.Value = '''  +0.31%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E49')
$cell.Value = '''  +4.87%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D50')
$cell.Value = '''20.68'
$cell.Style = 'Normal'
$cell = $ws.Range('E50')
$cell.Value = '''  -0.59%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D51')
$cell.Value = '''6.10'
$cell.Style = 'Normal'
$cell = $ws.Range('E51')
$cell.Value = '''  +1.98%  '
$cell.Style = 'Normal'
